$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric-looking accuracy figure into a cell as TEXT
# (matches the workbook's existing convention of storing these figures
# as shared strings, e.g. " 0.98130 "), while preserving the cell's
# original "centered / General" style.
function Set-TextValue($addr, $text) {
    $rng = $ws.Range($addr)
    $rng.Value = "'" + $text
    $rng.ClearFormats()
    $rng.HorizontalAlignment = -4108
}

# --- Header (No attack -> No attack (3 tries)) ---
$ws.Range("A2").Value = "No attack (3 tries)"

# --- kernel_size row: add a value of 3 for the softRmax and cons_softmax columns ---
$ws.Range("E11").Value = 3
$ws.Range("H11").Value = 3

# --- Per-epoch test accuracy labels (drop trailing ": ") ---
$ws.Range("A13").Value = "[epoch 1], test_accuracy"
$ws.Range("A14").Value = "[epoch 2], test_accuracy"
$ws.Range("A15").Value = "[epoch 3], test_accuracy"
$ws.Range("A16").Value = "[epoch 4], test_accuracy"
$ws.Range("A17").Value = "[epoch 5], test_accuracy"

# --- Updated accuracy figures: softmax (B col) ---
Set-TextValue "B13" " 0.97740"
Set-TextValue "B14" " 0.97870"
Set-TextValue "B15" " 0.98350"
Set-TextValue "B16" " 0.97240"
Set-TextValue "B17" " 0.97610"
Set-TextValue "B19" " 0.9835"

# --- Updated accuracy figures: softRmax (E col) ---
Set-TextValue "E13" " 0.98510"
Set-TextValue "E14" " 0.98950"
Set-TextValue "E15" " 0.98920"
Set-TextValue "E16" " 0.99060"
Set-TextValue "E17" " 0.99200"
Set-TextValue "E19" " 0.992"

# --- Updated accuracy figures: cons_softmax (H col) ---
Set-TextValue "H13" " 0.97820"
Set-TextValue "H14" " 0.98530"
Set-TextValue "H15" " 0.98950"
Set-TextValue "H16" " 0.98890"
Set-TextValue "H17" " 0.98830"
Set-TextValue "H19" " 0.9895"

# --- Selection moved to H19 ---
$ws.Range("H19").Select()
